$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the formatting of the
# existing header cells (bold, centered, bordered - same style as H1).
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J for rows 2-41.
$data = @(
    @(6, 7),
    @(8, 8),
    @(8, 9),
    @(6, 8),
    @(9, 9),
    @(8, 8),
    @(8, 9),
    @(6, 7),
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 9),
    @(8, 9),
    @(8, 9),
    @(5, 7),
    @(9, 9),
    @(6, 8),
    @(7, 8),
    @(5, 7),
    @(5, 5),
    @(4, 6),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(7, 8),
    @(9, 9),
    @(4, 6),
    @(6, 7),
    @(5, 6),
    @(4, 6),
    @(4, 5),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(3, 3),
    @(2, 2),
    @(5, 5),
    @(9, 9)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Range("I$row").Value = $data[$idx][0]
    $ws.Range("J$row").Value = $data[$idx][1]
}
